$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 31756
$ws.Range("J40").Value = 28159.6
$ws.Range("L40").Value = 28159.6
$ws.Range("N40").Value = -28509.6

$ws.Range("H137").Value = 5306
$ws.Range("I137").Value = 2997.8333
$ws.Range("J137").Value = 7284.4287
$ws.Range("K137").Value = 8993.499899999999
$ws.Range("L137").Value = 21853.2861
$ws.Range("M137").Value = -6443.499899999999
$ws.Range("N137").Value = -26953.2861

$ws.Range("H138").Value = 7567.943
$ws.Range("I138").Value = 2748.25
$ws.Range("J138").Value = 8996
$ws.Range("K138").Value = 8244.75
$ws.Range("L138").Value = 26988
$ws.Range("M138").Value = -3104.75
$ws.Range("N138").Value = -37268

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4803
$ws.Range("I32").Value = 4106.926
$ws.Range("J32").Value = 8561.799999999999
$ws.Range("K32").Value = 4106.926
$ws.Range("L32").Value = 8561.799999999999
$ws.Range("M32").Value = -3819.926
$ws.Range("N32").Value = -9135.799999999999

$ws.Range("H45").Value = 2575.0527
$ws.Range("I45").Value = 1309.9231
$ws.Range("J45").Value = 5316.1665
$ws.Range("K45").Value = 1309.9231
$ws.Range("L45").Value = 5316.1665
$ws.Range("M45").Value = -932.9231
$ws.Range("N45").Value = -6070.1665

$ws.Range("H74").Value = 5051.4873
$ws.Range("I74").Value = 1968.3334
$ws.Range("J74").Value = 5612.0605
$ws.Range("K74").Value = 1968.3334
$ws.Range("L74").Value = 5612.0605
$ws.Range("M74").Value = -1094.3334
$ws.Range("N74").Value = -7360.0605

$ws.Range("H77").Value = 5051.4873
$ws.Range("I77").Value = 1968.3334
$ws.Range("J77").Value = 5612.0605
$ws.Range("K77").Value = 9841.666999999999
$ws.Range("L77").Value = 28060.3025
$ws.Range("M77").Value = -5473.666999999999
$ws.Range("N77").Value = -36796.3025

$ws.Range("H102").Value = 1246709.1
$ws.Range("I102").Value = 1713644.1
$ws.Range("J102").Value = 1549.3334
$ws.Range("K102").Value = 1713644.1
$ws.Range("L102").Value = 1549.3334
$ws.Range("M102").Value = -1712022.1
$ws.Range("N102").Value = -4793.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 44282.11
$ws.Range("I58").Value = 29329.5
$ws.Range("K58").Value = 29329.5
$ws.Range("M58").Value = -29035.5

$ws.Range("H86").Value = 2547.2856
$ws.Range("I86").Value = 2520.8
$ws.Range("J86").Value = 2613.5
$ws.Range("K86").Value = 2520.8
$ws.Range("L86").Value = 2613.5
$ws.Range("M86").Value = -1397.8
$ws.Range("N86").Value = -4859.5

$ws.Range("H89").Value = 2547.2856
$ws.Range("I89").Value = 2520.8
$ws.Range("J89").Value = 2613.5
$ws.Range("K89").Value = 12604
$ws.Range("L89").Value = 13067.5
$ws.Range("M89").Value = -6988
$ws.Range("N89").Value = -24299.5

$ws.Range("H105").Value = 533.6
$ws.Range("I105").Value = 489
$ws.Range("K105").Value = 489
$ws.Range("M105").Value = 1258

$ws.Range("H134").Value = 3528.1064
$ws.Range("I134").Value = 2294.9143
$ws.Range("J134").Value = 7124.9165
$ws.Range("K134").Value = 6884.742899999999
$ws.Range("L134").Value = 21374.7495
$ws.Range("M134").Value = -4349.742899999999
$ws.Range("N134").Value = -26444.7495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15387863
$ws.Range("J31").Value = 5320.853
$ws.Range("L31").Value = 5320.853
$ws.Range("N31").Value = -5910.853

$ws.Range("H34").Value = 15387863
$ws.Range("J34").Value = 5320.853
$ws.Range("L34").Value = 5320.853
$ws.Range("N34").Value = -5724.853

$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("M86").Value = -2877

$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("M89").Value = -14384

$ws.Range("H134").Value = 2959.4
$ws.Range("I134").Value = 2999
$ws.Range("J134").Value = 2933
$ws.Range("K134").Value = 8997
$ws.Range("L134").Value = 8799
$ws.Range("N134").Value = -13869
$ws.Range("M134").Value = -6462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 797.6667
$ws.Range("I5").Value = 576.0714
$ws.Range("J5").Value = 1036.3077
$ws.Range("K5").Value = 1728.2142
$ws.Range("L5").Value = 3108.9231
$ws.Range("M5").Value = -1616.2142
$ws.Range("N5").Value = -3332.9231

$ws.Range("H68").Value = 339918.47
$ws.Range("J68").Value = 433840.38
$ws.Range("L68").Value = 1301521.14
$ws.Range("N68").Value = -1303143.14

$ws.Range("H71").Value = 339918.47
$ws.Range("J71").Value = 433840.38
$ws.Range("L71").Value = 3904563.42
$ws.Range("N71").Value = -3912675.42

$ws.Range("H86").Value = 1249.5
$ws.Range("I86").Value = 1333
$ws.Range("J86").Value = 999
$ws.Range("K86").Value = 3999
$ws.Range("L86").Value = 2997
$ws.Range("M86").Value = -2813
$ws.Range("N86").Value = -5369

$ws.Range("H89").Value = 1249.5
$ws.Range("I89").Value = 1333
$ws.Range("J89").Value = 999
$ws.Range("K89").Value = 11997
$ws.Range("L89").Value = 8991
$ws.Range("M89").Value = -6069
$ws.Range("N89").Value = -20847

$ws.Range("H107").Value = 1324786.4
$ws.Range("J107").Value = 1986644.6
$ws.Range("L107").Value = 5959933.800000001
$ws.Range("N107").Value = -5963773.800000001

$ws.Range("H131").Value = 1572768.4
$ws.Range("J131").Value = 2072573.6
$ws.Range("L131").Value = 6217720.800000001
$ws.Range("N131").Value = -6227800.800000001

$ws.Range("H135").Value = 797.6667
$ws.Range("I135").Value = 576.0714
$ws.Range("J135").Value = 1036.3077
$ws.Range("K135").Value = 5184.6426
$ws.Range("L135").Value = 9326.7693
$ws.Range("M135").Value = -2649.6426
$ws.Range("N135").Value = -14396.7693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 160.42857
$ws.Range("I2").Value = 159
$ws.Range("J2").Value = 162.33333
$ws.Range("K2").Value = 159
$ws.Range("L2").Value = 162.33333
$ws.Range("M2").Value = -46
$ws.Range("N2").Value = -388.33333

$ws.Range("H70").Value = 2273997.8
$ws.Range("I70").Value = 3407017.5
$ws.Range("J70").Value = 7958.143
$ws.Range("K70").Value = 3407017.5
$ws.Range("L70").Value = 7958.143
$ws.Range("M70").Value = -3406747.5
$ws.Range("N70").Value = -8498.143

$ws.Range("H73").Value = 2273997.8
$ws.Range("I73").Value = 3407017.5
$ws.Range("J73").Value = 7958.143
$ws.Range("K73").Value = 3407017.5
$ws.Range("L73").Value = 7958.143
$ws.Range("M73").Value = -3406081.5
$ws.Range("N73").Value = -9830.143

$ws.Range("H102").Value = 392717.3
$ws.Range("I102").Value = 562203.1
$ws.Range("K102").Value = 562203.1
$ws.Range("M102").Value = -560581.1

$ws.Range("H122").Value = 1107499.4
$ws.Range("I122").Value = 1577856.2
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 4733568.6
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -4731118.6
$ws.Range("N122").Value = -34900

$ws.Range("H126").Value = 3960.7827
$ws.Range("I126").Value = 2146.1538
$ws.Range("J126").Value = 6319.8
$ws.Range("K126").Value = 6438.4614
$ws.Range("L126").Value = 18959.4
$ws.Range("M126").Value = -3968.4614
$ws.Range("N126").Value = -23899.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1734.091
$ws.Range("J22").Value = 2330
$ws.Range("L22").Value = 2330
$ws.Range("N22").Value = -2920

$ws.Range("H27").Value = 1734.091
$ws.Range("J27").Value = 2330
$ws.Range("L27").Value = 2330
$ws.Range("N27").Value = -2544

$ws.Range("H40").Value = 20005798
$ws.Range("I40").Value = 50002500
$ws.Range("K40").Value = 50002500
$ws.Range("M40").Value = -50002364

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H29").Value = 11
$ws.Range("J29").Value = 11
$ws.Range("L29").Value = 11
$ws.Range("N29").Value = -591

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H81").Value = 1753464
$ws.Range("I81").Value = 2609265.8
$ws.Range("K81").Value = 5218531.6
$ws.Range("M81").Value = -5217470.6

$ws.Range("H84").Value = 1753464
$ws.Range("I84").Value = 2609265.8
$ws.Range("K84").Value = 26092658
$ws.Range("M84").Value = -26087354

$ws.Range("H122").Value = 3607.718
$ws.Range("I122").Value = 2709.4062
$ws.Range("J122").Value = 7714.2856
$ws.Range("K122").Value = 8128.2186
$ws.Range("L122").Value = 23142.8568
$ws.Range("M122").Value = -5678.2186
$ws.Range("N122").Value = -28042.8568

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
